$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B4 cell text from "Var3 " to lowercase "var3 " (case sensitivity fix)
$ws.Range("B4").Value = "var3 "

# Update the active selection on the sheet (E19 -> E14)
$ws.Range("E14").Select()
